$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.251.35'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.558.85'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = "'288.71"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = "'0.3806"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = "'0.3311"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = "'44.53"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.82%  '
$ws.Range('D10').Value = "'1.143"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').Value = "'0.07398"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = "'20.29"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.14%  '
$ws.Range('D14').Value = "'5.858"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('D15').Value = "'6.766"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.45%  '
$ws.Range('D16').Value = '1.549.08'
$ws.Range('E16').Value = '  -1.96%  '
$ws.Range('D17').Value = "'0.00001078"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.48%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = "'86.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = "'0.06647"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').Value = "'6.417"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = "'16.18"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('D23').Value = "'11.75"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '22.236.45'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').Value = "'2.284"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.87%  '
$ws.Range('D26').Value = "'2.573"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'151.23"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('D28').Value = "'19.30"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('D29').Value = "'4.941"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('D30').Value = "'123.25"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').Value = '1.729.43'
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').Value = "'1.088"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('D33').Value = "'5.932"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.82%  '
$ws.Range('D34').Value = "'1.923"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').Value = "'9.402"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').Value = "'0.08228"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').Value = "'0.02348"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.06321"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = "'5.343"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('D40').Value = "'0.2161"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('D41').Value = "'1.238"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.82%  '
$ws.Range('D42').Value = "'11.05"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('D43').Value = "'0.6086"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').Value = "'1.000"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = "'13.77"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.81%  '
$ws.Range('D46').Value = "'3.747"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = "'0.5896"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').Value = "'122.84"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('D49').Value = "'1.975"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.16%  '
$ws.Range('D50').Value = "'1.180"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('D51').Value = "'0.07072"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.84%  '
